# Register & log in TCs 03.03
# Fix a couple of misspelled Polish place names (missing "ó" diacritics)
# in the test-data sheet, and move the active cell selection to H5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Street address: "plac Podgorna 949" -> "plac Podgórna 949"
$ws.Range("G2").Value = "plac Podgórna 949"

# City: "Dąbrowa Gornicza" -> "Dąbrowa Górnicza"
$ws.Range("H2").Value = "Dąbrowa Górnicza"

# Update the sheet's active cell / selection to H5
$ws.Range("H5").Select()
